# Sync attendance_reports: normalize "Recorded By" (column G) values so that
# the "System" token is listed first instead of last, e.g.
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"         -> "System, backup@backdoor.com"
#   "system, backup@backdoor.com, System" -> "System, backup@backdoor.com, system"
#
# The rule observed in the data: whenever the comma-separated "Recorded By"
# text ends with ", System" the first and last comma-separated tokens are
# swapped (so "System" moves from the end to the front, and whatever token
# used to be first moves to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($null -ne $value -and $value -is [string] -and $value.EndsWith(", System")) {
        $parts = $value -split ", "
        if ($parts.Length -ge 2) {
            $first = $parts[0]
            $last = $parts[$parts.Length - 1]
            $parts[0] = $last
            $parts[$parts.Length - 1] = $first
            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
